# Shell scripts for model experiments
# Renames the "ADAM_*" model entries to lowercase "adam_*" naming, renames
# the idx header, and records which shell script ("sh") was used to launch
# each experiment in a new "Status"-adjacent column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: Model_idx -> idx ---------------------------------------
$ws.Range("A1").Value = "idx"

# --- Model_Name (column B) renames -----------------------------------------
$ws.Range("B4").Value  = "adam_wd=1e-2"
$ws.Range("B5").Value  = "adam_wd=1e-3"
$ws.Range("B6").Value  = "adam_dp=2e-1"
$ws.Range("B7").Value  = "adam_wd=1e-2_dp=2e-1"
$ws.Range("B8").Value  = "adam_wd=1e-3_dp=2e-1"
$ws.Range("B9").Value  = "adam_dp=4e-1"
$ws.Range("B10").Value = "adam_wd=1e-2_dp=4e-1"
$ws.Range("B11").Value = "adam_wd=1e-3_dp=4e-1"
$ws.Range("B12").Value = "adam_imgaug"
$ws.Range("B13").Value = "adam_wd=1e-2_imgaug"
$ws.Range("B14").Value = "adam_wd_1e-3_imgaug"
$ws.Range("B15").Value = "adam_dp=2e-1_imgaug"
$ws.Range("B16").Value = "adam_wd=1e-2_dp=2e-1_imgaug"
$ws.Range("B17").Value = "adam_wd=1e-3_dp=2e-1_imgaug"
$ws.Range("B18").Value = "adam_dp=4e-1_imgaug"
$ws.Range("B19").Value = "adam_wd=1e-2_dp=4e-1_imgaug"
$ws.Range("B20").Value = "adam_wd=1e-3_dp=4e-1_imgaug"

# --- New column G: which shell script launched the run ---------------------
$ws.Range("G4").Value  = "sh "
$ws.Range("G5").Value  = "sh "
$ws.Range("G6").Value  = "sh "
$ws.Range("G7").Value  = "sh - on HPC"
$ws.Range("G8").Value  = "sh"
$ws.Range("G9").Value  = "sh"
$ws.Range("G10").Value = "sh"
$ws.Range("G11").Value = "sh"
$ws.Range("G12").Value = "sh"
$ws.Range("G13").Value = "sh "
$ws.Range("G14").Value = "sh"
$ws.Range("G15").Value = "sh"
$ws.Range("G16").Value = "sh"
$ws.Range("G17").Value = "sh"
$ws.Range("G18").Value = "sh"
$ws.Range("G19").Value = "sh"
$ws.Range("G20").Value = "sh"

# --- Column width tweaks (idx column narrower, Model_Name / new G column
#     widened to fit their longer content) -----------------------------------
$ws.Columns.Item(1).ColumnWidth = 5.833333333333333
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(7).ColumnWidth = 9.5

# --- Restore the selection to where the editor left off ---------------------
[void]$ws.Range("G21").Select()
